$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.057.59'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.301.87'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.86%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.62'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.30%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.20'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +5.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.507'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.72%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.44'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.36%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.53%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '49.30'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +3.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.118'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.87'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +17.09%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.87%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.663.73'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.02%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.340.26'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.56%  '
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +4.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.944.85'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.07%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.32'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +8.25%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.49%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.81'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '236.72'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +12.88%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.47'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.78%  '
$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.52'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.30%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '168.04'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.09'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -8.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.79'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.73%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.57'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.11%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.97'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0700'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.103'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.57%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.82'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.110'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -2.29%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.002.79'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.81%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.08'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '17.58'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.87'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +2.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '55.03'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +5.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.529.19'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.77%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.27%  '
